$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 432-433; everything from old row 432 downward
# shifts down by two rows (old 432 -> new 434, ..., old 495 -> new 497).
$ws.Rows("432:433").Insert()

# New row 432 (Fruta, Terminal Hortofrutícola Agro Chillán - Pera)
$ws.Range("A432").Value2 = 7
$ws.Range("B432").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C432").Value2 = "Ñuble"
$ws.Range("D432").Value2 = 45154
$ws.Range("E432").Value2 = 16
$ws.Range("F432").Value2 = "Fruta"
$ws.Range("G432").Value2 = 100104
$ws.Range("H432").Value2 = "Frutos de pepita"
$ws.Range("I432").Value2 = 100104005
$ws.Range("J432").Value2 = "Pera"
$ws.Range("K432").Value2 = "Forelle"
$ws.Range("L432").Value2 = "Especial"
$ws.Range("M432").Value2 = 80
$ws.Range("N432").Value2 = 14000
$ws.Range("O432").Value2 = 14000
$ws.Range("P432").Value2 = 14000
$ws.Range("Q432").Value2 = "`$/bandeja 18 kilos granel"
$ws.Range("R432").Value2 = "Región de O'Higgins"
$ws.Range("S432").Value2 = 778
$ws.Range("T432").Value2 = 18

# New row 433
$ws.Range("A433").Value2 = 7
$ws.Range("B433").Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C433").Value2 = "Ñuble"
$ws.Range("D433").Value2 = 45154
$ws.Range("E433").Value2 = 16
$ws.Range("F433").Value2 = "Fruta"
$ws.Range("G433").Value2 = 100104
$ws.Range("H433").Value2 = "Frutos de pepita"
$ws.Range("I433").Value2 = 100104005
$ws.Range("J433").Value2 = "Pera"
$ws.Range("K433").Value2 = "Forelle"
$ws.Range("L433").Value2 = "Primera"
$ws.Range("M433").Value2 = 80
$ws.Range("N433").Value2 = 12000
$ws.Range("O433").Value2 = 12000
$ws.Range("P433").Value2 = 12000
$ws.Range("Q433").Value2 = "`$/bandeja 18 kilos granel"
$ws.Range("R433").Value2 = "Región de O'Higgins"
$ws.Range("S433").Value2 = 667
$ws.Range("T433").Value2 = 18
